$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-12-08 Sunday"; new = "2024-12-09 Monday"},
    @{old = "163÷7="; new = "674÷3="},
    @{old = "910÷3="; new = "218÷9="},
    @{old = "915÷3="; new = "188÷7="},
    @{old = "434÷2="; new = "426÷4="},
    @{old = "755÷5="; new = "624÷7="},
    @{old = "185÷9="; new = "312÷9="},
    @{old = "341÷9="; new = "116÷4="},
    @{old = "563÷6="; new = "438÷3="},
    @{old = "158÷4="; new = "595÷3="},
    @{old = "392÷3="; new = "644÷5="},
    @{old = "552÷2="; new = "720÷4="},
    @{old = "705÷7="; new = "569÷6="},
    @{old = "431÷4="; new = "313÷5="},
    @{old = "523÷6="; new = "986÷5="},
    @{old = "963÷3="; new = "714÷9="},
    @{old = "174÷5="; new = "874÷2="},
    @{old = "482÷9="; new = "945÷5="},
    @{old = "753÷4="; new = "637÷9="},
    @{old = "763÷4="; new = "793÷3="},
    @{old = "690÷3="; new = "156÷6="},
    @{old = "297÷2="; new = "373÷8="},
    @{old = "952÷7="; new = "110÷2="},
    @{old = "774÷6="; new = "159÷4="},
    @{old = "320÷7="; new = "764÷7="},
    @{old = "710÷9="; new = "881÷5="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2) | Out-Null
}
